$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1503
$ws1.Range("F5").Value = 208
$ws1.Range("F6").Value = 45
$ws1.Range("F7").Value = 175
$ws1.Range("F8").Value = 9876
$ws1.Range("F14").Value = 6830

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1503
$ws4.Range("F5").Value = 208
$ws4.Range("F7").Value = 45
$ws4.Range("F8").Value = 175
$ws4.Range("F11").Value = 9876
$ws4.Range("F17").Value = 6830
